# Update cryptos sheet with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.147.61'
$ws.Range('E2').Value = '  +2.35%  '
$ws.Range('D3').Value = '3.186.17'
$ws.Range('E3').Value = '  +1.89%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '532.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.80'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.526'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +10.69%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.30'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.439'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.113'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.75%  '
$ws.Range('E12').Value = '  +2.11%  '
$ws.Range('D13').Value = '3.722.62'
$ws.Range('E13').Value = '  +1.70%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.97'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000172'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.64%  '
$ws.Range('D16').Value = '59.114.72'
$ws.Range('E16').Value = '  +2.18%  '
$ws.Range('D17').Value = '3.216.53'
$ws.Range('E17').Value = '  +2.98%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.26'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.98%  '
$ws.Range('E19').Value = '  +2.37%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.20'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '375.52'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.91%  '
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.533'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.88'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.92%  '
$ws.Range('E25').Value = '  +0.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.48'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +16.04%  '
$ws.Range('D28').Value = '0.0₃0872'
$ws.Range('E28').Value = '  +0.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '22.44'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.70%  '
$ws.Range('E30').Value = '  +0.88%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.05'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.52%  '
$ws.Range('E32').Value = '  +0.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.17'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.36'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.51%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '157.41'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.33%  '
$ws.Range('E36').Value = '  +3.39%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0714'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.38%  '
$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.37'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.37%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '2.726.44'
$ws.Range('E39').Value = '  +7.18%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.71'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.82%  '
$ws.Range('E41').Value = '  +5.19%  '
$ws.Range('E42').Value = '  +8.81%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.726'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.88%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '39.20'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.89%  '
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('D46').Value = '3.224.10'
$ws.Range('E46').Value = '  +1.73%  '
$ws.Range('E47').Value = '  +11.78%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.988'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.99%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.23'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.33'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.78%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.758'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.10%  '
